# Apply crypto price/volume updates from the GitHub Actions scrape run.
# $q is a literal apostrophe used to force Excel to keep numeric-looking
# strings (e.g. "233.93") stored as text instead of auto-converting them.
$q = "'"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.772.69'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").Value = '2.087.52'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = $q + '233.93'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D8").Value = $q + '58.24'
$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").Value = $q + '0.0782'
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("E11").Value = '  +2.80%  '

$ws.Range("D12").Value = $q + '15.19'
$ws.Range("E12").Value = '  +1.63%  '

$ws.Range("D13").Value = '2.396.11'
$ws.Range("E13").Value = '  +0.41%  '

$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").Value = $q + '5.35'
$ws.Range("E16").Value = '  +0.96%  '

$ws.Range("D17").Value = '2.087.22'
$ws.Range("E17").Value = '  +0.58%  '

$ws.Range("D18").Value = '37.733.25'
$ws.Range("E18").Value = '  +0.33%  '

$ws.Range("D19").Value = $q + '6.11'
$ws.Range("E19").Value = '  -1.06%  '

$ws.Range("E20").Value = '  -0.18%  '

$ws.Range("E21").Value = '  +0.56%  '

$ws.Range("D22").Value = $q + '229.74'
$ws.Range("E22").Value = '  +0.75%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("E24").Value = '  -0.87%  '

$ws.Range("D25").Value = $q + '2.40'
$ws.Range("E25").Value = '  +0.17%  '

$ws.Range("E26").Value = '  +8.13%  '

$ws.Range("D27").Value = $q + '171.36'
$ws.Range("E27").Value = '  +1.35%  '

$ws.Range("D28").Value = $q + '0.134'
$ws.Range("E28").Value = '  -3.53%  '

$ws.Range("D29").Value = $q + '19.53'
$ws.Range("E29").Value = '  +0.42%  '

$ws.Range("E30").Value = '  -0.69%  '

$ws.Range("E31").Value = '  +0.57%  '

$ws.Range("D32").Value = $q + '4.69'
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").Value = $q + '0.0636'
$ws.Range("E33").Value = '  +1.05%  '

$ws.Range("D34").Value = $q + '4.62'
$ws.Range("E34").Value = '  -0.59%  '

$ws.Range("E35").Value = '  +0.81%  '

$ws.Range("E36").Value = '  -0.33%  '

$ws.Range("E37").Value = '  -1.81%  '

$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("D40").Value = $q + '0.0237'
$ws.Range("E40").Value = '  +10.11%  '

$ws.Range("D41").Value = $q + '101.19'
$ws.Range("E41").Value = '  +3.14%  '

$ws.Range("E42").Value = '  -0.67%  '

$ws.Range("E43").Value = '  +4.46%  '

$ws.Range("E44").Value = '  +1.53%  '

$ws.Range("D45").Value = $q + '16.67'
$ws.Range("E45").Value = '  +0.37%  '

$ws.Range("D46").Value = '1.452.14'
$ws.Range("E46").Value = '  +0.05%  '

$ws.Range("E47").Value = '  -4.15%  '

$ws.Range("E48").Value = '  -0.26%  '

$ws.Range("D49").Value = $q + '7.20'

$ws.Range("E50").Value = '  -1.95%  '

$ws.Range("D51").Value = '2.280.24'
$ws.Range("E51").Value = '  +0.40%  '
